$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-27 20:18:42'
$ws.Range('E3').Value = '2026-02-27 20:18:45'
$ws.Range('H3').Value = '39%'
$ws.Range('N3').Value = '0.5 °C 19:59 TU'
$ws.Range('O3').Value = '4.6 °C'
$ws.Range('E4').Value = '2026-02-27 20:18:47'
$ws.Range('O4').Value = '9.8 °C'
$ws.Range('E5').Value = '2026-02-27 20:18:50'
$ws.Range('H5').Value = '41%'
$ws.Range('N5').Value = '0.8 °C 19:43 TU'
$ws.Range('O5').Value = '5.0 °C'
$ws.Range('E6').Value = '2026-02-27 20:18:53'
$ws.Range('J6').Value = '1024.3 hPa'
$ws.Range('E7').Value = '2026-02-27 20:18:56'
$ws.Range('E8').Value = '2026-02-27 20:18:58'
$ws.Range('H8').Value = '64%'
$ws.Range('N8').Value = '8.2 °C 19:58 TU'
$ws.Range('O8').Value = '11.9 °C'
$ws.Range('E9').Value = '2026-02-27 20:19:00'
$ws.Range('E10').Value = '2026-02-27 20:19:03'
$ws.Range('K10').Value = '14.2 MJ/m2'
$ws.Range('O10').Value = '11.1 °C'
$ws.Range('E11').Value = '2026-02-27 20:19:05'
$ws.Range('E12').Value = '2026-02-27 20:19:08'
$ws.Range('E13').Value = '2026-02-27 20:19:11'
$ws.Range('H13').Value = '65%'
$ws.Range('E14').Value = '2026-02-27 20:19:13'
$ws.Range('O14').Value = '10.8 °C'
$ws.Range('E15').Value = '2026-02-27 20:19:16'
$ws.Range('E16').Value = '2026-02-27 20:19:18'
$ws.Range('H16').Value = '43%'
$ws.Range('E17').Value = '2026-02-27 20:19:21'
$ws.Range('E18').Value = '2026-02-27 20:19:23'
$ws.Range('H18').Value = '81%'
$ws.Range('J18').Value = '1024.5 hPa'
$ws.Range('E19').Value = '2026-02-27 20:19:25'
$ws.Range('E20').Value = '2026-02-27 20:19:28'
$ws.Range('N20').Value = '-0.2 °C 19:59 TU'
$ws.Range('O20').Value = '3.3 °C'
$ws.Range('E21').Value = '2026-02-27 20:19:31'
$ws.Range('J21').Value = '1024.2 hPa'
$ws.Range('K21').Value = '15.4 MJ/m2'
$ws.Range('E22').Value = '2026-02-27 20:19:33'
$ws.Range('E23').Value = '2026-02-27 20:19:36'
$ws.Range('H23').Value = '41%'
$ws.Range('E24').Value = '2026-02-27 20:19:39'
$ws.Range('E25').Value = '2026-02-27 20:19:42'
$ws.Range('O25').Value = '6.2 °C'
$ws.Range('E26').Value = '2026-02-27 20:19:45'
$ws.Range('H26').Value = '45%'
$ws.Range('O26').Value = '10.4 °C'
$ws.Range('E27').Value = '2026-02-27 20:19:48'
$ws.Range('N27').Value = '2.7 °C 19:57 TU'
$ws.Range('O27').Value = '5.7 °C'
$ws.Range('E28').Value = '2026-02-27 20:19:50'
$ws.Range('E29').Value = '2026-02-27 20:19:53'
$ws.Range('E30').Value = '2026-02-27 20:19:56'
$ws.Range('H30').Value = '91%'
$ws.Range('E31').Value = '2026-02-27 20:19:59'
$ws.Range('H31').Value = '92%'
$ws.Range('J31').Value = '1024.1 hPa'
$ws.Range('E32').Value = '2026-02-27 20:20:01'
$ws.Range('K32').Value = '14.2 MJ/m2'
$ws.Range('E33').Value = '2026-02-27 20:20:04'
$ws.Range('E34').Value = '2026-02-27 20:20:07'
$ws.Range('O34').Value = '4.7 °C'
$ws.Range('E35').Value = '2026-02-27 20:20:10'
$ws.Range('J35').Value = '1022.5 hPa'
$ws.Range('E36').Value = '2026-02-27 20:20:12'
$ws.Range('E37').Value = '2026-02-27 20:20:15'
$ws.Range('J37').Value = '1024.8 hPa'
$ws.Range('O37').Value = '8.5 °C'
$ws.Range('E38').Value = '2026-02-27 20:20:18'
$ws.Range('E39').Value = '2026-02-27 20:20:21'
$ws.Range('E40').Value = '2026-02-27 20:20:23'
$ws.Range('J40').Value = '1024.7 hPa'
$ws.Range('O40').Value = '9.2 °C'
$ws.Range('E41').Value = '2026-02-27 20:20:26'
$ws.Range('J41').Value = '1024.6 hPa'
$ws.Range('O41').Value = '11.1 °C'
$ws.Range('E42').Value = '2026-02-27 20:20:29'
$ws.Range('E43').Value = '2026-02-27 20:20:31'
$ws.Range('E44').Value = '2026-02-27 20:20:34'
$ws.Range('H44').Value = '59%'
$ws.Range('E45').Value = '2026-02-27 20:20:37'
$ws.Range('O45').Value = '11.9 °C'
$ws.Range('E46').Value = '2026-02-27 20:20:40'
$ws.Range('J46').Value = '1024.1 hPa'
